# Insert a new data row before the current row 639 (shifts existing rows
# 639..665 down to 640..666) and populate it with a new price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("639:639").Insert()

$ws.Range("A639").Value = 11
$ws.Range("B639").Value = "Vega Monumental Concepción"
$ws.Range("C639").Value = "Bíobío"
$ws.Range("D639").Value = 44939
$ws.Range("E639").Value = 8
$ws.Range("F639").Value = "Fruta"
$ws.Range("G639").Value = 100108
$ws.Range("H639").Value = "Tropicales y subtropicales"
$ws.Range("I639").Value = 100108006
$ws.Range("J639").Value = "Plátano"
$ws.Range("K639").Value = "Sin especificar"
$ws.Range("L639").Value = "Pintón"
$ws.Range("M639").Value = 1200
$ws.Range("N639").Value = 21000
$ws.Range("O639").Value = 22000
$ws.Range("P639").Value = 21500
$ws.Range("Q639").Value = "$/caja 20 kilos"
$ws.Range("R639").Value = "Ecuador"
$ws.Range("S639").Value = 1075
$ws.Range("T639").Value = 20
